# Manual renamed to User Guide. Comments of review incorporated.
# Add "Additional Effort [h]" of 1 hour to the 2013-07-16 entry (row 34)
# and log a new day of work: 2013-07-18, 1.5h, "Revision of manual" (row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# Row 34 (2013-07-16): record additional effort of 1 hour in column C
$ws.Cells.Item(34, 3).Value = 1

# New row 36 (2013-07-18): 1.5h "Revision of manual"
$ws.Cells.Item(36, 1).Value = Get-Date -Year 2013 -Month 7 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(36, 2).Value = 1.5
$ws.Cells.Item(36, 4).Value = "Revision of manual"

# Move the active selection past the last used row, as after a fresh data entry
$ws.Range("A37").Select()
